$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First table (rows 1-7): label addition ---
$ws.Range("G9").Value = "Processes"

# --- Second table body (rows 10-15): newly measured data ---
$ws.Range("B10").Value = 300
$ws.Range("C10").Value = 1050
$ws.Range("D10").Value = 1160
$ws.Range("E10").Value = 442

$ws.Range("B15").Value = 2715

# --- New "combined scope" results table starting at row 18 ---
$ws.Range("I18").Value = "min"
$ws.Range("J18").Value = "avg"

# Row 19 - fully literal values, no shared formulas
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = 12
$ws.Range("C19").Value = 24
$ws.Range("A19").Font.Bold = $true
$ws.Range("B19").Font.Bold = $true
$ws.Range("C19").Font.Bold = $true

$ws.Range("D19").Value = 2715
$ws.Range("E19").Value = 3110
$ws.Range("F19").Value = 2690
$ws.Range("G19").Value = 3220
$ws.Range("H19").Value = 4440
$ws.Range("I19").Formula = "=MIN(D19:H19)"
$ws.Range("J19").Formula = "=AVERAGE(D19:H19)"

# Rows 20-30 - A/B/C labels (bold), D-H raw data, I/J shared MIN/AVERAGE formulas
$rows = @(20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30)
$aVals = @(2, 2, 2, 4, 2, 4, 5, 2, 4, 6, 8)
$bVals = @(10, 8, 6, 6, 4, 4, 4, 2, 2, 2, 2)
$cVals = @(20, 16, 12, 24, 8, 16, 9, 4, 8, 12, 16)

for ($idx = 0; $idx -lt $rows.Count; $idx++) {
    $r = $rows[$idx]
    $ws.Cells.Item($r, 1).Value = $aVals[$idx]
    $ws.Cells.Item($r, 2).Value = $bVals[$idx]
    $ws.Cells.Item($r, 3).Value = $cVals[$idx]
    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 2).Font.Bold = $true
    $ws.Cells.Item($r, 3).Font.Bold = $true
}

# Row 20 has two computed-average formula cells among its raw values
$ws.Range("D20").Formula = "=(3200+2940)/2"
$ws.Range("E20").Value = 2740
$ws.Range("F20").Value = 2920
$ws.Range("G20").Value = 2395
$ws.Range("H20").Formula = "=(3090+3290)/2"

$ws.Range("D21").Value = 1350
$ws.Range("E21").Value = 1580
$ws.Range("F21").Value = 2920
$ws.Range("G21").Value = 1975
$ws.Range("H21").Value = 1645

$ws.Range("D22").Value = 500
$ws.Range("E22").Value = 485
$ws.Range("F22").Value = 455
$ws.Range("G22").Value = 480
$ws.Range("H22").Value = 510

$ws.Range("D23").Value = 1470
$ws.Range("E23").Value = 1525
$ws.Range("F23").Value = 1230
$ws.Range("G23").Value = 1395

$ws.Range("D24").Value = 330
$ws.Range("E24").Value = 370
$ws.Range("F24").Value = 255
$ws.Range("G24").Value = 355
$ws.Range("H24").Value = 360

$ws.Range("D25").Value = 407
$ws.Range("E25").Value = 450
$ws.Range("F25").Value = 500
$ws.Range("G25").Value = 560
$ws.Range("H25").Value = 510

$ws.Range("D26").Value = 750
$ws.Range("E26").Value = 790
$ws.Range("F26").Value = 680
$ws.Range("G26").Value = 860
$ws.Range("H26").Value = 940

$ws.Range("D27").Value = 350
$ws.Range("E27").Value = 360
$ws.Range("F27").Value = 250
$ws.Range("G27").Value = 240
$ws.Range("H27").Value = 225

$ws.Range("D28").Value = 390
$ws.Range("E28").Value = 310
$ws.Range("F28").Value = 840
$ws.Range("G28").Value = 280
$ws.Range("H28").Value = 300

$ws.Range("D29").Value = 490
$ws.Range("E29").Value = 1020
$ws.Range("F29").Value = 380
$ws.Range("G29").Value = 390
$ws.Range("H29").Value = 480

$ws.Range("D30").Value = 430
$ws.Range("E30").Value = 460
$ws.Range("F30").Value = 610
$ws.Range("G30").Value = 450
$ws.Range("H30").Value = 490

# Shared MIN/AVERAGE formulas for rows 20-30
$ws.Range("I20:I30").Formula = "=MIN(D20:H20)"
$ws.Range("J20:J30").Formula = "=AVERAGE(D20:H20)"

$excel.Calculate()

# --- View state to match the saved selection ---
$ws.Range("K18:L18").Select()
